# UART Round 1 Changes
# - Typography sheet: clear the "Wildcard Characters" (H) value on the
#   "Default" (row 4) and "SingleUseId28"/0x... (row 17) typography rows.
# - Translation sheet: the translation table is trimmed down to just the
#   three still-relevant rows (SingleUseId8, SingleUseId20, SingleUseId21),
#   with row 6 (SingleUseId21) now pointing at a brand new "CAN Message
#   Output" text instead of the old "Current Time" text. All the other
#   previously-present rows (7-34) are cleared out.

$wb = $excel.ActiveWorkbook

# ---- Typography sheet --------------------------------------------------
$tj = $wb.Worksheets.Item("Typography")
$tj.Range("H4").ClearContents()
$tj.Range("H17").ClearContents()

# ---- Translation sheet --------------------------------------------------
$tr = $wb.Worksheets.Item("Translation")

# Clear every existing data row (4 through 34) first ...
$tr.Range("B4:F34").ClearContents()

# ... then write back only the three rows that remain.
$tr.Range("B4").Value = "SingleUseId8"
$tr.Range("C4").Value = "Typography_01"
$tr.Range("D4").Value = "Center"
$tr.Range("E4").Value = "LTR"
$tr.Range("F4").Value = "Cancel"

$tr.Range("B5").Value = "SingleUseId20"
$tr.Range("C5").Value = "Typography_06"
$tr.Range("D5").Value = "Center"
$tr.Range("E5").Value = "LTR"
$tr.Range("F5").Value = "STOP MOTOR"

$tr.Range("B6").Value = "SingleUseId21"
$tr.Range("C6").Value = "Default"
$tr.Range("D6").Value = "Left"
$tr.Range("E6").Value = "LTR"
$tr.Range("F6").Value = "CAN Message Output"
